$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------
# Add relationships for the "Checker" model: mark challengeId columns as
# foreign keys (challengeId *), note the FK constraint names and index,
# and add the "* foreign keys" footnote. New shared strings must be
# created in this exact order so the rebuilt sharedStrings table lines
# up with the target workbook:
#   1) "challengeId *"
#   2) "Tasks_ChallengeId_fkey"
#   3) "to Challenges"
#   4) "ChallengeId_foreign_idx"
#   5) "* foreign keys"
# -------------------------------------------------------------------

# 1) Users table (row 4): F4 challengeId -> challengeId *  (creates new string #1)
$ws.Range("F4").Value = "challengeId *"

# 2) Task table row 11: challengeId -> challengeId * (reuses string #1),
#    then add the FK constraint name (creates new string #2)
$ws.Range("E11").Value = "challengeId *"
$ws.Range("F11").Value = "Tasks_ChallengeId_fkey"

# 3) Users table row 4: add "to Challenges" annotation cell (creates new string #3)
#    New cell needs the same style as the other label cells (style index 1)
$ws.Range("B3").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").Value = "to Challenges"

# Task table row 11: add matching "to Challenges" annotation (reuses string #3)
$ws.Range("B3").Copy()
$ws.Range("H11").PasteSpecial(-4122)
$ws.Range("H11").Value = "to Challenges"

# 4) Users table row 4: add index name cell (creates new string #4)
$ws.Range("G4").Value = "ChallengeId_foreign_idx"

# 5) New footnote row 18 (creates new string #5)
$ws.Range("G18").Value = "* foreign keys"
$ws.Rows.Item(18).RowHeight = 15.75

# -------------------------------------------------------------------
# Remaining shared-string reindex fixups (string "challengeId" was
# removed from the table, shifting every later index down by one):
# -------------------------------------------------------------------

# Challenge table title (row 6)
$ws.Range("B6").Value = "Challenge"

# Challenge table header row 7: (createdAt) / (updatedat)
$ws.Range("H7").Value = "(createdAt)"
$ws.Range("I7").Value = "(updatedat)"

# Task table title (row 10)
$ws.Range("B10").Value = "Task"

# Task table row 12: type label
$ws.Range("E12").Value = "interger"

# Checker table title (row 15)
$ws.Range("B15").Value = "Checker"

# Checker table header row 16
$ws.Range("C16").Value = "taskId"
$ws.Range("D16").Value = "(createdAt)"
$ws.Range("E16").Value = "(updatedat)"

# -------------------------------------------------------------------
# View state (best-effort)
# -------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G20").Select()
